# Weekly update: insert a new price record at the top of the data block
# (row 196), pushing the existing rows 196:264 down to 197:265.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(196).Insert()

$ws.Range("A196").Value = 7
$ws.Range("B196").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C196").Value = "Ñuble"
$ws.Range("D196").Value = 44627
$ws.Range("E196").Value = 16
$ws.Range("F196").Value = 100114013
$ws.Range("G196").Value = "Zanahoria"
$ws.Range("H196").Value = "Sin especificar"
$ws.Range("I196").Value = "Primera"
$ws.Range("J196").Value = 100
$ws.Range("K196").Value = 7000
$ws.Range("L196").Value = 7500
$ws.Range("M196").Value = 7250
$ws.Range("N196").Value = "$/saco 20 kilos"
$ws.Range("O196").Value = "Provincia de Diguillín"
$ws.Range("P196").Value = 362
$ws.Range("Q196").Value = 20
$ws.Range("R196").Value = "Hortaliza"
